# Weekly update: insert a new price record at the top of the data table
# (row 11, just below the header + existing row 10) and push the older
# records down by one row. This mirrors the upstream "Fruta / hortaliza,
# semanal" commit, which prepends the newest observation and keeps the
# rest of the history intact (the oldest visible row is pushed from 21
# to the newly created row 22).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Push rows 11..21 down to 12..22, duplicating formatting (incl. the
# date style on column D) the same way Excel's own "Insert Row" does.
$ws.Rows.Item(11).Insert()

# Populate the freshly inserted row 11 with the new weekly record.
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(11, 3).Value = "La Araucanía"
$ws.Cells.Item(11, 4).Value = 44438
$ws.Cells.Item(11, 5).Value = 9
$ws.Cells.Item(11, 6).Value = 100114002
$ws.Cells.Item(11, 7).Value = "Camote"
$ws.Cells.Item(11, 8).Value = "Sin especificar"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 40
$ws.Cells.Item(11, 11).Value = 20000
$ws.Cells.Item(11, 12).Value = 20000
$ws.Cells.Item(11, 13).Value = 20000
$ws.Cells.Item(11, 14).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(11, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(11, 16).Value = 1333
$ws.Cells.Item(11, 17).Value = 15
$ws.Cells.Item(11, 18).Value = "Hortaliza"
